# Laboratórios.xlsx refactor:
#   - row 2 -> "Unidade Interlagos" branch data (was "Unidade Pinheiros")
#   - row 3 -> "Unidade Primavera" branch data (was "Unidade Penha")
#   - E2 becomes the numeric building number 6 (was 123)
#   - E3 becomes the text address suffix "S/N" instead of the number 212
#   - header/data rows grow to 15pt tall
#   - view scrolled back so column A is visible (topLeftCell -> A1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: "Unidade Interlagos" ---
$ws.Range("A2").Value = "Unidade Interlagos"
$ws.Range("B2").Value = "São Paulo"
$ws.Range("C2").Value = "Interlagos"
$ws.Range("D2").Value = "Sabará"
$ws.Range("E2").Value = 6

# --- Row 3: "Unidade Primavera" ---
$ws.Range("A3").Value = "Unidade Primavera"
$ws.Range("B3").Value = "São Paulo"
$ws.Range("C3").Value = "Primavera Interlagos"
$ws.Range("D3").Value = "Rua 10"
$ws.Range("E3").Value = "S/N"

# --- Row heights: header + both data rows now 15pt ---
$ws.Range("A1:E3").RowHeight = 15

# --- Scroll the view back to the left so column A is the top-left cell ---
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
